# Applies the crypto price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.582.45'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '3.695.29'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''676.41'
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").Value = '''160.71'
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.67%  '
$ws.Range("E9").Value = '  +1.32%  '
$ws.Range("D10").Value = '''7.11'
$ws.Range("E10").Value = '  +0.12%  '
$ws.Range("E11").Value = '  +1.36%  '
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("D13").Value = '''32.55'
$ws.Range("E13").Value = '  +0.59%  '
$ws.Range("D14").Value = '3.699.20'
$ws.Range("E14").Value = '  +0.27%  '
$ws.Range("D15").Value = '69.580.06'
$ws.Range("E15").Value = '  +0.34%  '
$ws.Range("E16").Value = '  +2.03%  '
$ws.Range("E17").Value = '  +1.38%  '
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("E19").Value = '  +0.43%  '
$ws.Range("D20").Value = '''9.81'
$ws.Range("E20").Value = '  -2.58%  '
$ws.Range("E21").Value = '  +0.87%  '
$ws.Range("D22").Value = '''80.61'
$ws.Range("E22").Value = '  +1.26%  '
$ws.Range("D23").Value = '3.843.47'
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("E25").Value = '  +3.86%  '
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("E28").Value = '  +0.73%  '
$ws.Range("E29").Value = '  +1.88%  '
$ws.Range("D30").Value = '''2.01'
$ws.Range("E30").Value = '  +0.45%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = '''6.59'
$ws.Range("E31").Value = '  +0.42%  '
$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").Value = '''1.00'
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("D33").Value = '''26.94'
$ws.Range("E33").Value = '  +0.94%  '
$ws.Range("D34").Value = '3.685.74'
$ws.Range("E34").Value = '  +0.82%  '
$ws.Range("D35").Value = '''0.162'
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("E36").Value = '  +4.29%  '
$ws.Range("E37").Value = '  +1.88%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").Value = '''1.00'
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '''2.23'
$ws.Range("E40").Value = '  -1.76%  '
$ws.Range("D41").Value = '''0.0902'
$ws.Range("E41").Value = '  +0.48%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Value = '''0.944'
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = '''166.82'
$ws.Range("E43").Value = '  +1.42%  '
$ws.Range("D44").Value = '''47.01'
$ws.Range("E44").Value = '  -1.14%  '
$ws.Range("E45").Value = '  +2.43%  '
$ws.Range("D46").Value = '''28.11'
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("E47").Value = '  +1.16%  '
$ws.Range("E48").Value = '  +0.60%  '
$ws.Range("E49").Value = '  -2.59%  '
$ws.Range("D50").Value = '''7.86'
$ws.Range("E50").Value = '  +0.93%  '
$ws.Range("E51").Value = '  +2.05%  '
